$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.138.14"
$ws.Range("E2").Value = "'  +0.27%  "
$ws.Range("D3").Value = "'2.327.30"
$ws.Range("E3").Value = "'  -0.69%  "
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'528.04"
$ws.Range("E5").Value = "'  +1.97%  "
$ws.Range("D6").Value = "'132.43"
$ws.Range("E6").Value = "'  -2.41%  "
$ws.Range("E7").Value = "'  -0.41%  "
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "'  -0.60%  "
$ws.Range("D9").Value = "'2.355.32"
$ws.Range("E9").Value = "'  -0.10%  "
$ws.Range("E10").Value = "'  -1.38%  "
$ws.Range("E11").Value = "'  +0.48%  "
$ws.Range("D12").Value = "'5.31"
$ws.Range("E12").Value = "'  -2.19%  "
$ws.Range("D13").Value = "'0.345"
$ws.Range("E13").Value = "'  +0.46%  "
$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'2.749.04"
$ws.Range("E14").Value = "'  -0.49%  "
$ws.Range("B15").Value = "'Avalanche"
$ws.Range("C15").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'23.55"
$ws.Range("E15").Value = "'  -1.80%  "
$ws.Range("D16").Value = "'57.169.18"
$ws.Range("E16").Value = "'  +0.34%  "
$ws.Range("E17").Value = "'  -1.28%  "
$ws.Range("D18").Value = "'2.348.41"
$ws.Range("E18").Value = "'  -0.11%  "
$ws.Range("D19").Value = "'336.82"
$ws.Range("D20").Value = "'10.45"
$ws.Range("E20").Value = "'  -1.27%  "
$ws.Range("D21").Value = "'6.91"
$ws.Range("E21").Value = "'  +2.08%  "
$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "'  -1.39%  "
$ws.Range("E23").Value = "'  -0.07%  "
$ws.Range("D24").Value = "'61.74"
$ws.Range("E24").Value = "'  +0.85%  "
$ws.Range("D25").Value = "'8.91"
$ws.Range("E25").Value = "'  +11.72%  "
$ws.Range("E26").Value = "'  +0.12%  "
$ws.Range("D27").Value = "'0.991"
$ws.Range("E27").Value = "'  -0.54%  "
$ws.Range("D28").Value = "'1.33"
$ws.Range("E28").Value = "'  +2.96%  "
$ws.Range("D29").Value = "'169.93"
$ws.Range("E29").Value = "'  -0.27%  "
$ws.Range("D30").Value = "'1.71"
$ws.Range("E30").Value = "'  +0.69%  "
$ws.Range("E31").Value = "'  -2.58%  "
$ws.Range("D32").Value = "'6.14"
$ws.Range("E32").Value = "'  -1.70%  "
$ws.Range("D33").Value = "'18.54"
$ws.Range("E33").Value = "'  -0.17%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "'  -0.05%  "
$ws.Range("D35").Value = "'0.992"
$ws.Range("E35").Value = "'  -0.38%  "
$ws.Range("D36").Value = "'1.26"
$ws.Range("E36").Value = "'  -0.69%  "
$ws.Range("B37").Value = "'NEARProtocol"
$ws.Range("C37").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.02"
$ws.Range("E37").Value = "'  -0.08%  "
$ws.Range("B38").Value = "'SuiNetwork"
$ws.Range("C38").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'0.909"
$ws.Range("E38").Value = "'  -1.05%  "
$ws.Range("E39").Value = "'  +1.26%  "
$ws.Range("D40").Value = "'38.95"
$ws.Range("E40").Value = "'  +1.43%  "
$ws.Range("D41").Value = "'148.43"
$ws.Range("E41").Value = "'  +0.59%  "
$ws.Range("D42").Value = "'0.379"
$ws.Range("E42").Value = "'  -1.37%  "
$ws.Range("B43").Value = "'Bittensor"
$ws.Range("C43").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'286.15"
$ws.Range("E43").Value = "'  +1.65%  "
$ws.Range("B44").Value = "'Filecoin"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'3.60"
$ws.Range("E44").Value = "'  -0.98%  "
$ws.Range("E45").Value = "'  -2.57%  "
$ws.Range("E46").Value = "'  -0.48%  "
$ws.Range("E47").Value = "'  -0.61%  "
$ws.Range("E48").Value = "'  -0.57%  "
$ws.Range("D49").Value = "'18.71"
$ws.Range("E49").Value = "'  +3.72%  "
$ws.Range("E50").Value = "'  -0.88%  "
$ws.Range("B51").Value = "'Polygon"
$ws.Range("C51").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D51").Value = "'0.379"
$ws.Range("E51").Value = "'  -0.86%  "
